$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws "D2" "35.063.22"
Set-TextValue $ws "E2" "  +0.26%  "
Set-TextValue $ws "D3" "1.820.88"
Set-TextValue $ws "E3" "  +0.02%  "
Set-TextValue $ws "E4" "  -0.20%  "
Set-TextValue $ws "D5" "230.76"
Set-TextValue $ws "E5" "  -0.31%  "
Set-TextValue $ws "D6" "0.618"
Set-TextValue $ws "E6" "  +0.74%  "
Set-TextValue $ws "E7" "  -0.26%  "
Set-TextValue $ws "D8" "40.00"
Set-TextValue $ws "E8" "  -4.73%  "
Set-TextValue $ws "E9" "  +4.86%  "
Set-TextValue $ws "D10" "0.0684"
Set-TextValue $ws "E10" "  -0.01%  "
Set-TextValue $ws "D11" "0.0991"
Set-TextValue $ws "E11" "  -1.14%  "
Set-TextValue $ws "D12" "2.082.71"
Set-TextValue $ws "E12" "  -0.21%  "
Set-TextValue $ws "D13" "11.34"
Set-TextValue $ws "E13" "  +2.33%  "
Set-TextValue $ws "D14" "0.668"
Set-TextValue $ws "E14" "  +1.34%  "
Set-TextValue $ws "D15" "1.811.66"
Set-TextValue $ws "E15" "  -0.30%  "
Set-TextValue $ws "D16" "4.63"
Set-TextValue $ws "E16" "  -0.55%  "
Set-TextValue $ws "D17" "35.005.11"
Set-TextValue $ws "E17" "  +0.05%  "
Set-TextValue $ws "D18" "69.71"
Set-TextValue $ws "E18" "  +0.22%  "
Set-TextValue $ws "D19" "0.0₃0788"
Set-TextValue $ws "E19" "  +0.22%  "
Set-TextValue $ws "D20" "241.11"
Set-TextValue $ws "E20" "  +1.23%  "
Set-TextValue $ws "D21" "12.08"
Set-TextValue $ws "E21" "  +2.76%  "
Set-TextValue $ws "D22" "4.68"
Set-TextValue $ws "E22" "  +2.39%  "
Set-TextValue $ws "E23" "  -0.03%  "
Set-TextValue $ws "D25" "173.55"
Set-TextValue $ws "E25" "  +0.73%  "
Set-TextValue $ws "E26" "  +1.15%  "
Set-TextValue $ws "E27" "  +3.25%  "
Set-TextValue $ws "D28" "17.38"
Set-TextValue $ws "E28" "  -0.16%  "
Set-TextValue $ws "E29" "  -4.24%  "
Set-TextValue $ws "E30" "  -0.08%  "
Set-TextValue $ws "E31" "  +2.73%  "
Set-TextValue $ws "E32" "  +0.16%  "
Set-TextValue $ws "D33" "3.96"
Set-TextValue $ws "E33" "  -0.13%  "
Set-TextValue $ws "E34" "  +12.97%  "
Set-TextValue $ws "E35" "  +3.21%  "
Set-TextValue $ws "D36" "0.694"
Set-TextValue $ws "E36" "  +2.98%  "
Set-TextValue $ws "D37" "92.89"
Set-TextValue $ws "E37" "  +0.32%  "
Set-TextValue $ws "E38" "  +5.92%  "
Set-TextValue $ws "D39" "1.340.48"
Set-TextValue $ws "E39" "  +2.19%  "
Set-TextValue $ws "E40" "  +0.87%  "
Set-TextValue $ws "D41" "0.988"
Set-TextValue $ws "E41" "  +0.15%  "
Set-TextValue $ws "D42" "14.68"
Set-TextValue $ws "E42" "  +0.10%  "
Set-TextValue $ws "E43" "  -0.94%  "
Set-TextValue $ws "E44" "  -1.28%  "
Set-TextValue $ws "E45" "  -0.82%  "
Set-TextValue $ws "E46" "  +2.20%  "
Set-TextValue $ws "D47" "6.22"
Set-TextValue $ws "E47" "  +0.74%  "
Set-TextValue $ws "D48" "1.999.01"
Set-TextValue $ws "E48" "  -0.04%  "
Set-TextValue $ws "E49" "  -0.11%  "
Set-TextValue $ws "E50" "  +3.78%  "
Set-TextValue $ws "D51" "97.08"
Set-TextValue $ws "E51" "  -2.81%  "
